# Add team record (Wins/Losses/Ties) columns to the data sheet, matching
# the existing header formatting (bold, centered, thin border).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): new columns AD, AE, AF.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

foreach ($col in @("AD1", "AE1", "AF1")) {
    $cell = $ws.Range($col)
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108   # xlCenter
    $cell.VerticalAlignment = -4160     # xlTop
    $cell.Borders.LineStyle = 1         # xlContinuous (thin box border)
}

# Data rows (2 through 51): every row gets the same team record (89-73-0).
for ($r = 2; $r -le 51; $r++) {
    $ws.Cells.Item($r, 30).Value = 89   # AD -> Wins
    $ws.Cells.Item($r, 31).Value = 73   # AE -> Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF -> Ties
}
